$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.227.24"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.76%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.055.16"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.44%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.27"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.34%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.616"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.36%  "

# Row 7
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.76"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +4.17%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.382"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.65%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.20"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.64%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0758"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.82%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.361.59"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.54%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.35"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.20%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.67"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.03%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.770"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.89%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.14"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.50%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.057.76"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.48%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "37.185.52"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.67%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.28"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +12.73%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.01"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.74%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0810"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.12%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "224.86"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.45%  "

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.00%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.43"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.90%  "

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.29%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.10"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.91%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.45"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +6.60%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.77"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.13%  "

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.85%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.02"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.23%  "

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.44%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.46"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.54%  "

# Row 34
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0613"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.65%  "

# Row 35
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.53"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.14%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.52"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.98%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.19%  "

# Row 38
$ws.Range("B38").Value = "THORChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.83"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.06%  "

# Row 39
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.27"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.81%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.74"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.56%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.63"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +13.98%  "

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.33%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.482.78"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.37%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "96.68"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.07%  "

# Row 45
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.16"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.51%  "

# Row 46
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0926"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.83%  "

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.89%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.35"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.41%  "

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.48%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.18"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.20%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.96"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.13%  "
